$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("B11").Value = 6
$ws.Range("B12").Value = 7
$ws.Range("B13").Value = 8
$ws.Range("B14").Value = 9
$ws.Range("B15").Value = 10
$ws.Range("B16").Value = 11
$ws.Range("B17").Value = 12
$ws.Range("B18").Value = 13
$ws.Range("B19").Value = 15
$ws.Range("B20").Value = 14
$ws.Range("B21").Value = 16
$ws.Range("B22").Value = 17
$ws.Range("B23").Value = 18
$ws.Range("B24").Value = 19
$ws.Range("B25").Value = 20
$ws.Range("B26").Value = 21
$ws.Range("B27").Value = 22
$ws.Range("B28").Value = 23
$ws.Range("B29").Value = 24
$ws.Range("B30").Value = 25
$ws.Range("B31").Value = 26
$ws.Range("B32").Value = 27
$ws.Range("B33").Value = 28
$ws.Range("B34").Value = 29
$ws.Range("B35").Value = 30
$ws.Range("B36").Value = 31
$ws.Range("B37").Value = 32
$ws.Range("B38").Value = 33
$ws.Range("B39").Value = 34
$ws.Range("B40").Value = 35
$ws.Range("B41").Value = 36
$ws.Range("B42").Value = 37
$ws.Range("B43").Value = 38
$ws.Range("B44").Value = 39
$ws.Range("B45").Value = 40
$ws.Range("B46").Value = 41
$ws.Range("B47").Value = 42
$ws.Range("B48").Value = 42
$ws.Range("B49").Value = 43
$ws.Range("B50").Value = 44
$ws.Range("B51").Value = 45
$ws.Range("B52").Value = 46
$ws.Range("B53").Value = 47
$ws.Range("B54").Value = 48
$ws.Range("B55").Value = 49
$ws.Range("B56").Value = 50
$ws.Range("B57").Value = 51
$ws.Range("B58").Value = 104
$ws.Range("B59").Value = 52
$ws.Range("B60").Value = 52
$ws.Range("B61").Value = 53
$ws.Range("B62").Value = 54
$ws.Range("B63").Value = 55
$ws.Range("B64").Value = 56
$ws.Range("B65").Value = 58
$ws.Range("B66").Value = 57
$ws.Range("B67").Value = 59
$ws.Range("B68").Value = 60
$ws.Range("B69").Value = 60
$ws.Range("B70").Value = 61
$ws.Range("B71").Value = 62
$ws.Range("B72").Value = 63
$ws.Range("B73").Value = 64
$ws.Range("B74").Value = 64
$ws.Range("B75").Value = 65
$ws.Range("B76").Value = 66
$ws.Range("B77").Value = 67
$ws.Range("B78").Value = 68
$ws.Range("B79").Value = 69
$ws.Range("B80").Value = 70
$ws.Range("B81").Value = 71
$ws.Range("B82").Value = 72
$ws.Range("B83").Value = 73
$ws.Range("B84").Value = 74
$ws.Range("B85").Value = 75
$ws.Range("B86").Value = 76
$ws.Range("B87").Value = 77
$ws.Range("B88").Value = 78
$ws.Range("B89").Value = 79
$ws.Range("B90").Value = 80
$ws.Range("B91").Value = 81
$ws.Range("B92").Value = 82
$ws.Range("B93").Value = 83
$ws.Range("B94").Value = 84
$ws.Range("B95").Value = 85
$ws.Range("B96").Value = 86
$ws.Range("B97").Value = 87
$ws.Range("B98").Value = 88
$ws.Range("B99").Value = 89
$ws.Range("B100").Value = 90
$ws.Range("B101").Value = 91
$ws.Range("B102").Value = 92
$ws.Range("B103").Value = 93
$ws.Range("B104").Value = 94
$ws.Range("B105").Value = 95
$ws.Range("B106").Value = 96
$ws.Range("B107").Value = 97
$ws.Range("B108").Value = 98
$ws.Range("B109").Value = 99
$ws.Range("B110").Value = 100
$ws.Range("B111").Value = 101
$ws.Range("B112").Value = 102
$ws.Range("B113").Value = 103

$ws.Activate()
$ws.Range("B59:B113").Select()
